# Apply balance-analysis fix: update computed Xcg/Ycg/Zcg values.
$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Range("C4").Value = 8.545043666543112
$wsGlobal.Range("C5").Value = 27.50726038674152
$wsGlobal.Range("C6").Value = 7.545759053209421
$wsGlobal.Range("C7").Value = 24.433376021409586
$wsGlobal.Range("C8").Value = 7.347756617256282
$wsGlobal.Range("C9").Value = 23.824303707023063

# --- POWER PLANT sheet (ENGINE 1 / ENGINE 2 balance) ---
$wsPower = $wb.Worksheets.Item("POWER PLANT")
$wsPower.Range("C13").Value = 33.40549999999999
$wsPower.Range("C15").Value = 4.437499999999999
$wsPower.Range("C17").Value = 1.6499999999999997
$wsPower.Range("C21").Value = 33.40549999999999
$wsPower.Range("C23").Value = -4.437499999999999
$wsPower.Range("C25").Value = 1.6499999999999997
